{"js": "const replacements = [\n  [\"328\u00f77=\", \"663\u00f78=\"],\n  [\"724\u00f75=\", \"669\u00f78=\"],\n  [\"677\u00f73=\", \"409\u00f78=\"],\n  [\"309\u00f78=\", \"559\u00f72=\"],\n  [\"418\u00f72=\", \"347\u00f74=\"],\n  [\"660\u00f72=\", \"148\u00f72=\"],\n  [\"300\u00f74=\", \"550\u00f79=\"],\n  [\"694\u00f75=\", \"710\u00f76=\"],\n  [\"978\u00f78=\", \"362\u00f74=\"],\n  [\"744\u00f79=\", \"885\u00f73=\"],\n  [\"562\u00f77=\", \"914\u00f74=\"],\n  [\"571\u00f77=\", \"470\u00f76=\"],\n  [\"613\u00f75=\", \"356\u00f73=\"],\n  [\"391\u00f78=\", \"939\u00f73=\"],\n  [\"797\u00f72=\", \"830\u00f72=\"],\n  [\"351\u00f79=\", \"303\u00f75=\"],\n  [\"829\u00f76=\", \"482\u00f73=\"],\n  [\"661\u00f74=\", \"324\u00f74=\"],\n  [\"888\u00f74=\", \"574\u00f78=\"],\n  [\"373\u00f78=\", \"316\u00f76=\"],\n  [\"471\u00f79=\", \"812\u00f78=\"],\n  [\"414\u00f75=\", \"405\u00f73=\"],\n  [\"795\u00f77=\", \"102\u00f79=\"],\n  [\"939\u00f72=\", \"687\u00f75=\"],\n  [\"487\u00f78=\", \"879\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldVal, newVal] of replacements) {\n  const results = body.search(oldVal, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldVal);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newVal, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"328\u00f77=\", \"663\u00f78=\"),\n    @(\"724\u00f75=\", \"669\u00f78=\"),\n    @(\"677\u00f73=\", \"409\u00f78=\"),\n    @(\"309\u00f78=\", \"559\u00f72=\"),\n    @(\"418\u00f72=\", \"347\u00f74=\"),\n    @(\"660\u00f72=\", \"148\u00f72=\"),\n    @(\"300\u00f74=\", \"550\u00f79=\"),\n    @(\"694\u00f75=\", \"710\u00f76=\"),\n    @(\"978\u00f78=\", \"362\u00f74=\"),\n    @(\"744\u00f79=\", \"885\u00f73=\"),\n    @(\"562\u00f77=\", \"914\u00f74=\"),\n    @(\"571\u00f77=\", \"470\u00f76=\"),\n    @(\"613\u00f75=\", \"356\u00f73=\"),\n    @(\"391\u00f78=\", \"939\u00f73=\"),\n    @(\"797\u00f72=\", \"830\u00f72=\"),\n    @(\"351\u00f79=\", \"303\u00f75=\"),\n    @(\"829\u00f76=\", \"482\u00f73=\"),\n    @(\"661\u00f74=\", \"324\u00f74=\"),\n    @(\"888\u00f74=\", \"574\u00f78=\"),\n    @(\"373\u00f78=\", \"316\u00f76=\"),\n    @(\"471\u00f79=\", \"812\u00f78=\"),\n    @(\"414\u00f75=\", \"405\u00f73=\"),\n    @(\"795\u00f77=\", \"102\u00f79=\"),\n    @(\"939\u00f72=\", \"687\u00f75=\"),\n    @(\"487\u00f78=\", \"879\u00f73=\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldVal = $pair[0]\n    $newVal = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute(\n        $oldVal,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        $wdFindContinue,\n        $false,\n        $newVal,\n        $wdReplaceAll\n    ) | Out-Null\n}\n"}
